$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Nid2"
$ws.Range("C2").Value = "Col13a1"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 52.269495
$ws.Range("H2").Value = 104.53899
$ws.Range("I2").Value = 0.4288758502612202
$ws.Range("J2").Value = 0.3474396247360552
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 0.4274105
$ws.Range("N2").Value = 0.854821
$ws.Range("O2").Value = 0.4422185568930271
$ws.Range("P2").Value = 0.3984644485961792
$ws.Range("Q2").Value = 22.3405309926975
$ws.Range("R2").Value = 89.36212397079001
$ws.Range("S2").Value = 0.1896568595887868
$ws.Range("T2").Value = 0.1384423384909156

# Row 3
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Nid2"
$ws.Range("C3").Value = "Col13a1"
$ws.Range("D3").Value = "Neutro"
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 52.269495
$ws.Range("H3").Value = 104.53899
$ws.Range("I3").Value = 0.4288758502612202
$ws.Range("J3").Value = 0.3474396247360552
$ws.Range("K3").Value = 2
$ws.Range("L3").Value = 0.6666666666666666
$ws.Range("M3").Value = 0.2122596666666667
$ws.Range("N3").Value = 0.636779
$ws.Range("O3").Value = 0.2196136114576696
$ws.Range("P3").Value = 0.2968268129966699
$ws.Range("Q3").Value = 11.094705585535
$ws.Range("R3").Value = 66.56823351321
$ws.Range("S3").Value = 0.09418697434284531
$ws.Range("T3").Value = 0.1031293965191622

# Row 4
$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Nid2"
$ws.Range("C4").Value = "Col13a1"
$ws.Range("D4").Value = "sCs"
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 52.269495
$ws.Range("H4").Value = 104.53899
$ws.Range("I4").Value = 0.4288758502612202
$ws.Range("J4").Value = 0.3474396247360552
$ws.Range("K4").Value = 2
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 0.326844
$ws.Range("N4").Value = 0.653688
$ws.Range("O4").Value = 0.3381678316493033
$ws.Range("P4").Value = 0.304708738407151
$ws.Range("Q4").Value = 17.08397082378
$ws.Range("R4").Value = 68.33588329512
$ws.Range("S4").Value = 0.1450320163295881
$ws.Range("T4").Value = 0.1058678897259773

# Row 5
$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Nid2"
$ws.Range("C5").Value = "Col13a1"
$ws.Range("D5").Value = "ECs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 27.97505033333333
$ws.Range("H5").Value = 83.925151
$ws.Range("I5").Value = 0.2295377733763979
$ws.Range("J5").Value = 0.2789286845927703
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 0.4274105
$ws.Range("N5").Value = 0.854821
$ws.Range("O5").Value = 0.4422185568930271
$ws.Range("P5").Value = 0.3984644485961792
$ws.Range("Q5").Value = 11.95683025049517
$ws.Range("R5").Value = 71.74098150297101
$ws.Range("S5").Value = 0.1015058628949494
$ws.Range("T5").Value = 0.1111431645039158

# Row 6
$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Nid2"
$ws.Range("C6").Value = "Col13a1"
$ws.Range("D6").Value = "Neutro"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 27.97505033333333
$ws.Range("H6").Value = 83.925151
$ws.Range("I6").Value = 0.2295377733763979
$ws.Range("J6").Value = 0.2789286845927703
$ws.Range("K6").Value = 2
$ws.Range("L6").Value = 0.6666666666666666
$ws.Range("M6").Value = 0.2122596666666667
$ws.Range("N6").Value = 0.636779
$ws.Range("O6").Value = 0.2196136114576696
$ws.Range("P6").Value = 0.2968268129966699
$ws.Range("Q6").Value = 5.937974858736555
$ws.Range("R6").Value = 53.44177372862899
$ws.Range("S6").Value = 0.05040961937714288
$ws.Range("T6").Value = 0.08279351250102536

# Row 7
$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Nid2"
$ws.Range("C7").Value = "Col13a1"
$ws.Range("D7").Value = "sCs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 27.97505033333333
$ws.Range("H7").Value = 83.925151
$ws.Range("I7").Value = 0.2295377733763979
$ws.Range("J7").Value = 0.2789286845927703
$ws.Range("K7").Value = 2
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 0.326844
$ws.Range("N7").Value = 0.653688
$ws.Range("O7").Value = 0.3381678316493033
$ws.Range("P7").Value = 0.304708738407151
$ws.Range("Q7").Value = 9.143477351148
$ws.Range("R7").Value = 54.860864106888
$ws.Range("S7").Value = 0.07762229110430567
$ws.Range("T7").Value = 0.08499200758782917

# Row 8
$ws.Range("A8").Value = "M1"
$ws.Range("B8").Value = "Nid2"
$ws.Range("C8").Value = "Col13a1"
$ws.Range("D8").Value = "ECs"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 0.7109406666666667
$ws.Range("H8").Value = 2.132822
$ws.Range("I8").Value = 0.005833331332203332
$ws.Range("J8").Value = 0.007088521472311936
$ws.Range("K8").Value = 2
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 0.4274105
$ws.Range("N8").Value = 0.854821
$ws.Range("O8").Value = 0.4422185568930271
$ws.Range("P8").Value = 0.3984644485961792
$ws.Range("Q8").Value = 0.3038635058103333
$ws.Range("R8").Value = 1.823181034862
$ws.Range("S8").Value = 0.002579607363605837
$ws.Range("T8").Value = 0.002824523799826952

# Row 9
$ws.Range("A9").Value = "M1"
$ws.Range("B9").Value = "Nid2"
$ws.Range("C9").Value = "Col13a1"
$ws.Range("D9").Value = "Neutro"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 0.7109406666666667
$ws.Range("H9").Value = 2.132822
$ws.Range("I9").Value = 0.005833331332203332
$ws.Range("J9").Value = 0.007088521472311936
$ws.Range("K9").Value = 2
$ws.Range("L9").Value = 0.6666666666666666
$ws.Range("M9").Value = 0.2122596666666667
$ws.Range("N9").Value = 0.636779
$ws.Range("O9").Value = 0.2196136114576696
$ws.Range("P9").Value = 0.2968268129966699
$ws.Range("Q9").Value = 0.1509040289264444
$ws.Range("R9").Value = 1.358136260338
$ws.Range("S9").Value = 0.001281078960694353
$ws.Range("T9").Value = 0.002104063237484814

# Row 10
$ws.Range("A10").Value = "M1"
$ws.Range("B10").Value = "Nid2"
$ws.Range("C10").Value = "Col13a1"
$ws.Range("D10").Value = "sCs"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 0.7109406666666667
$ws.Range("H10").Value = 2.132822
$ws.Range("I10").Value = 0.005833331332203332
$ws.Range("J10").Value = 0.007088521472311936
$ws.Range("K10").Value = 2
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 0.326844
$ws.Range("N10").Value = 0.653688
$ws.Range("O10").Value = 0.3381678316493033
$ws.Range("P10").Value = 0.304708738407151
$ws.Range("Q10").Value = 0.232366691256
$ws.Range("R10").Value = 1.394200147536
$ws.Range("S10").Value = 0.001972645007903143
$ws.Range("T10").Value = 0.00215993443500017

# Row 11
$ws.Range("A11").Value = "M2"
$ws.Range("B11").Value = "Nid2"
$ws.Range("C11").Value = "Col13a1"
$ws.Range("D11").Value = "ECs"
$ws.Range("E11").Value = 3
$ws.Range("F11").Value = 1
$ws.Range("G11").Value = 1.165253333333333
$ws.Range("H11").Value = 3.49576
$ws.Range("I11").Value = 0.009561007124768556
$ws.Range("J11").Value = 0.01161830186581401
$ws.Range("K11").Value = 2
$ws.Range("L11").Value = 1
$ws.Range("M11").Value = 0.4274105
$ws.Range("N11").Value = 0.854821
$ws.Range("O11").Value = 0.4422185568930271
$ws.Range("P11").Value = 0.3984644485961792
$ws.Range("Q11").Value = 0.4980415098266668
$ws.Range("R11").Value = 2.98824905896
$ws.Range("S11").Value = 0.004228054773159102
$ws.Range("T11").Value = 0.00462948024658554

# Row 12
$ws.Range("A12").Value = "M2"
$ws.Range("B12").Value = "Nid2"
$ws.Range("C12").Value = "Col13a1"
$ws.Range("D12").Value = "Neutro"
$ws.Range("E12").Value = 3
$ws.Range("F12").Value = 1
$ws.Range("G12").Value = 1.165253333333333
$ws.Range("H12").Value = 3.49576
$ws.Range("I12").Value = 0.009561007124768556
$ws.Range("J12").Value = 0.01161830186581401
$ws.Range("K12").Value = 2
$ws.Range("L12").Value = 0.6666666666666666
$ws.Range("M12").Value = 0.2122596666666667
$ws.Range("N12").Value = 0.636779
$ws.Range("O12").Value = 0.2196136114576696
$ws.Range("P12").Value = 0.2968268129966699
$ws.Range("Q12").Value = 0.2473362841155556
$ws.Range("R12").Value = 2.22602655704
$ws.Range("S12").Value = 0.002099727303842933
$ws.Range("T12").Value = 0.003448623515262837

# Row 13
$ws.Range("A13").Value = "M2"
$ws.Range("B13").Value = "Nid2"
$ws.Range("C13").Value = "Col13a1"
$ws.Range("D13").Value = "sCs"
$ws.Range("E13").Value = 3
$ws.Range("F13").Value = 1
$ws.Range("G13").Value = 1.165253333333333
$ws.Range("H13").Value = 3.49576
$ws.Range("I13").Value = 0.009561007124768556
$ws.Range("J13").Value = 0.01161830186581401
$ws.Range("K13").Value = 2
$ws.Range("L13").Value = 1
$ws.Range("M13").Value = 0.326844
$ws.Range("N13").Value = 0.653688
$ws.Range("O13").Value = 0.3381678316493033
$ws.Range("P13").Value = 0.304708738407151
$ws.Range("Q13").Value = 0.3808560604800001
$ws.Range("R13").Value = 2.28513636288
$ws.Range("S13").Value = 0.003233225047766523
$ws.Range("T13").Value = 0.003540198103965636

# Row 14
$ws.Range("A14").Value = "Neutro"
$ws.Range("B14").Value = "Nid2"
$ws.Range("C14").Value = "Col13a1"
$ws.Range("D14").Value = "ECs"
$ws.Range("E14").Value = 3
$ws.Range("F14").Value = 1
$ws.Range("G14").Value = 27.28148833333333
$ws.Range("H14").Value = 81.844465
$ws.Range("I14").Value = 0.22384703554817
$ws.Range("J14").Value = 0.2720134392567137
$ws.Range("K14").Value = 2
$ws.Range("L14").Value = 1
$ws.Range("M14").Value = 0.4274105
$ws.Range("N14").Value = 0.854821
$ws.Range("O14").Value = 0.4422185568930271
$ws.Range("P14").Value = 0.3984644485961792
$ws.Range("Q14").Value = 11.66039456929417
$ws.Range("R14").Value = 69.962367415765
$ws.Range("S14").Value = 0.09898931302489386
$ws.Range("T14").Value = 0.1083876850841767

# Row 15
$ws.Range("A15").Value = "Neutro"
$ws.Range("B15").Value = "Nid2"
$ws.Range("C15").Value = "Col13a1"
$ws.Range("D15").Value = "Neutro"
$ws.Range("E15").Value = 3
$ws.Range("F15").Value = 1
$ws.Range("G15").Value = 27.28148833333333
$ws.Range("H15").Value = 81.844465
$ws.Range("I15").Value = 0.22384703554817
$ws.Range("J15").Value = 0.2720134392567137
$ws.Range("K15").Value = 2
$ws.Range("L15").Value = 0.6666666666666666
$ws.Range("M15").Value = 0.2122596666666667
$ws.Range("N15").Value = 0.636779
$ws.Range("O15").Value = 0.2196136114576696
$ws.Range("P15").Value = 0.2968268129966699
$ws.Range("Q15").Value = 5.790759619803889
$ws.Range("R15").Value = 52.116836578235
$ws.Range("S15").Value = 0.04915985589082696
$ws.Range("T15").Value = 0.0807408822668336

# Row 16
$ws.Range("A16").Value = "Neutro"
$ws.Range("B16").Value = "Nid2"
$ws.Range("C16").Value = "Col13a1"
$ws.Range("D16").Value = "sCs"
$ws.Range("E16").Value = 3
$ws.Range("F16").Value = 1
$ws.Range("G16").Value = 27.28148833333333
$ws.Range("H16").Value = 81.844465
$ws.Range("I16").Value = 0.22384703554817
$ws.Range("J16").Value = 0.2720134392567137
$ws.Range("K16").Value = 2
$ws.Range("L16").Value = 1
$ws.Range("M16").Value = 0.326844
$ws.Range("N16").Value = 0.653688
$ws.Range("O16").Value = 0.3381678316493033
$ws.Range("P16").Value = 0.304708738407151
$ws.Range("Q16").Value = 8.91679077282
$ws.Range("R16").Value = 53.50074463692
$ws.Range("S16").Value = 0.07569786663244915
$ws.Range("T16").Value = 0.08288487190570343

# Row 17
$ws.Range("A17").Value = "sCs"
$ws.Range("B17").Value = "Nid2"
$ws.Range("C17").Value = "Col13a1"
$ws.Range("D17").Value = "ECs"
$ws.Range("E17").Value = 2
$ws.Range("F17").Value = 1
$ws.Range("G17").Value = 12.473357
$ws.Range("H17").Value = 24.946714
$ws.Range("I17").Value = 0.10234500235724
$ws.Range("J17").Value = 0.0829114280763349
$ws.Range("K17").Value = 2
$ws.Range("L17").Value = 1
$ws.Range("M17").Value = 0.4274105
$ws.Range("N17").Value = 0.854821
$ws.Range("O17").Value = 0.4422185568930271
$ws.Range("P17").Value = 0.3984644485961792
$ws.Range("Q17").Value = 5.331243752048501
$ws.Range("R17").Value = 21.324975008194
$ws.Range("S17").Value = 0.04525885924763212
$ws.Range("T17").Value = 0.03303725647075856

# Row 18
$ws.Range("A18").Value = "sCs"
$ws.Range("B18").Value = "Nid2"
$ws.Range("C18").Value = "Col13a1"
$ws.Range("D18").Value = "Neutro"
$ws.Range("E18").Value = 2
$ws.Range("F18").Value = 1
$ws.Range("G18").Value = 12.473357
$ws.Range("H18").Value = 24.946714
$ws.Range("I18").Value = 0.10234500235724
$ws.Range("J18").Value = 0.0829114280763349
$ws.Range("K18").Value = 2
$ws.Range("L18").Value = 0.6666666666666666
$ws.Range("M18").Value = 0.2122596666666667
$ws.Range("N18").Value = 0.636779
$ws.Range("O18").Value = 0.2196136114576696
$ws.Range("P18").Value = 0.2968268129966699
$ws.Range("Q18").Value = 2.647590599034333
$ws.Range("R18").Value = 15.885543594206
$ws.Range("S18").Value = 0.02247635558231718
$ws.Range("T18").Value = 0.02461033495690111

# Row 19
$ws.Range("A19").Value = "sCs"
$ws.Range("B19").Value = "Nid2"
$ws.Range("C19").Value = "Col13a1"
$ws.Range("D19").Value = "sCs"
$ws.Range("E19").Value = 2
$ws.Range("F19").Value = 1
$ws.Range("G19").Value = 12.473357
$ws.Range("H19").Value = 24.946714
$ws.Range("I19").Value = 0.10234500235724
$ws.Range("J19").Value = 0.0829114280763349
$ws.Range("K19").Value = 2
$ws.Range("L19").Value = 1
$ws.Range("M19").Value = 0.326844
$ws.Range("N19").Value = 0.653688
$ws.Range("O19").Value = 0.3381678316493033
$ws.Range("P19").Value = 0.304708738407151
$ws.Range("Q19").Value = 4.076841895308
$ws.Range("R19").Value = 16.307367581232
$ws.Range("S19").Value = 0.03460978752729068
$ws.Range("T19").Value = 0.02526383664867524
